$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update several "DATA EXTRACAO" timestamps to newer re-scrape values ---
$ws.Range("C7").Value  = "2025-02-04 09:52:03"   # barramansa.rj.gov.br
$ws.Range("C10").Value = "2025-02-04 10:13:54"   # buzios.rj.gov.br
$ws.Range("C12").Value = "2025-02-04 10:40:27"   # cachoeirasdemacacu.rj.gov.br
$ws.Range("C17").Value = "2025-02-04 10:18:33"   # carmo.rj.gov.br
$ws.Range("C21").Value = "2025-02-04 10:41:41"   # duasbarras.rj.gov.br
$ws.Range("C22").Value = "2025-02-04 10:42:49"   # duquedecaxias.rj.gov.br
$ws.Range("C35").Value = "2025-02-04 10:43:01"   # mage.rj.gov.br
$ws.Range("C48").Value = "2025-02-04 10:43:02"   # paracambi.rj.gov.br
$ws.Range("C82").Value = "2025-02-04 10:43:04"   # sjb.rj.gov.br

# --- Fix Guapimirim's "DATA EXTRACAO" (row 23): it was stored as plain text;
#     convert it to a real Excel date/time value (same instant,
#     2025-02-03 09:51:40) formatted as yyyy-mm-dd hh:mm:ss ---
$cell = $ws.Range("C23")
$cell.NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$cell.Value = Get-Date -Year 2025 -Month 2 -Day 3 -Hour 9 -Minute 51 -Second 40

# --- Move the active selection (cosmetic, matches the saved view state) ---
$ws.Range("B5").Select() | Out-Null
